# Update "想去人数" (F) / "最低票价" (G) counters with refreshed scrape values,
# matching the gh-pages regeneration at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 26165
$ws.Range("F6").Value = 247
$ws.Range("F7").Value = 586
$ws.Range("F9").Value = 488
$ws.Range("F12").Value = 215
$ws.Range("F15").Value = 286
$ws.Range("F17").Value = 362
$ws.Range("F18").Value = 51
$ws.Range("F20").Value = 170
$ws.Range("F21").Value = 19
$ws.Range("F22").Value = 422
$ws.Range("F23").Value = 97
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 141
$ws.Range("F8").Value = 107
$ws.Range("F9").Value = 107
$ws.Range("F10").Value = 430
$ws.Range("G17").Value = 880
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4956
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 4956
$ws.Range("F6").Value = 26165
$ws.Range("F9").Value = 247
$ws.Range("F11").Value = 586
$ws.Range("F15").Value = 141
$ws.Range("F16").Value = 141
$ws.Range("F18").Value = 107
$ws.Range("F19").Value = 107
$ws.Range("F20").Value = 430
$ws.Range("F21").Value = 489
$ws.Range("F25").Value = 215
$ws.Range("F29").Value = 286
$ws.Range("F33").Value = 362
$ws.Range("F34").Value = 51
$ws.Range("F37").Value = 170
$ws.Range("F39").Value = 19
$ws.Range("F40").Value = 422
$ws.Range("F41").Value = 97
$ws.Range("G43").Value = 880
